# repull data, push all data, mean calculation
# Update the dSF column (F) with repulled values for each row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    2  = 5
    3  = 3
    4  = -3
    6  = 1
    7  = 11
    8  = -7
    9  = -3
    10 = -7
    11 = -2
    12 = 2
    13 = 3
    14 = -4
    15 = 2
    16 = 0
    17 = 4
    18 = 3
    19 = -1
    20 = -2
    21 = 1
    23 = -1
    24 = 0
    25 = -2
    26 = -2
    27 = -5
    28 = -2
    29 = -4
    30 = -1
    31 = -5
    32 = -1
    33 = -1
    34 = -2
    37 = -2
}

foreach ($row in $values.Keys) {
    $ws.Range("F$row").Value = $values[$row]
}
